# Changes to image generation
#
# The per-topic "image prompt" descriptions on the image_prompts sheet are
# replaced by a single, unified art-director style prompt that is now used
# for every topic row (Blockchain, Sports, Culture, Exclusive, Opinion,
# Business, World News, General).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("image_prompts")

$bullet = [char]0x2022
$mid    = [char]0x00B7
$rq     = [char]0x201D
$lq     = [char]0x201C
$dash   = [char]0x2014
$apos   = [char]0x2019

$newPrompt = "You are an art director creating descriptive prompts for DALL${mid}E 3.`nYour task is to write a single, production-ready image prompt that transforms a short news article summary into a symbolic, hand-drawn sketch.`nFollow these creative guidelines:`n${bullet} The image should look like an **artist${apos}s sketch** ${dash} loose, imperfect, expressive line work on textured paper.`n${bullet} Emphasize **hand-drawn pencil or ink strokes**, light shading, and visible texture.`n${bullet} Avoid digital gloss or rendering; do not use vector or 3D styles.`n${bullet} Keep the color palette muted and minimal. Use the provided accent color sparingly, as a small highlight or ink tone.`n${bullet} Depict ideas **conceptually and symbolically**, not literally.`n${bullet} No text, logos, faces, or copyrighted symbols.`n${bullet} The overall style should resemble a **newspaper editorial illustration** ${dash} subtle, elegant, and human in feel.`nOutput format:`n1. One paragraph describing the scene in natural, directive language suitable for DALL${mid}E 3.`n2. Optionally, a few short style tags (e.g. ${lq}hand-drawn sketch, pencil on paper, minimal color, conceptual illustration${rq}).`nOutput only the prompt text${dash}no explanations or extra words."

# Apply the same prompt text to every topic row (rows 2-9); the header in
# row 1 ("desc_image_prompt") is left untouched.
$ws.Range("B2:B9").Value = $newPrompt

# The long text needs considerably taller rows once it wraps in column B.
$ws.Rows("2:9").RowHeight = 356

# Reflect the saved scroll position / selection from the workbook view.
$ws.Range("A8").Select()
$ws.Range("B1").Select()
